$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '62.048.30'
Set-TextValue 'E2' '  -0.09%  '
Set-TextValue 'D3' '2.411.51'
Set-TextValue 'E3' '  -0.32%  '
Set-TextValue 'E4' '  +0.05%  '
Set-TextValue 'D5' '563.77'
Set-TextValue 'E5' '  +1.63%  '
Set-TextValue 'D6' '142.71'
Set-TextValue 'E6' '  -0.25%  '
Set-TextValue 'E7' '  +0.08%  '
Set-TextValue 'D8' '0.529'
Set-TextValue 'E8' '  -0.66%  '
Set-TextValue 'E9' '  +0.64%  '
Set-TextValue 'E10' '  -1.75%  '
Set-TextValue 'E11' '  -1.74%  '
Set-TextValue 'D12' '0.351'
Set-TextValue 'E12' '  -0.55%  '
Set-TextValue 'D13' '25.63'
Set-TextValue 'E13' '  -2.39%  '
Set-TextValue 'E14' '  -0.67%  '
Set-TextValue 'D15' '2.847.02'
Set-TextValue 'E15' '  -0.27%  '
Set-TextValue 'D16' '61.953.43'
Set-TextValue 'E16' '  -0.08%  '
Set-TextValue 'D17' '2.409.90'
Set-TextValue 'E17' '  -0.34%  '
Set-TextValue 'E18' '  +1.43%  '
Set-TextValue 'D19' '6.87'
Set-TextValue 'E19' '  +1.81%  '
Set-TextValue 'D20' '321.76'
Set-TextValue 'E20' '  -0.86%  '
Set-TextValue 'E21' '  -1.11%  '
Set-TextValue 'E22' '  -0.09%  '
Set-TextValue 'D23' '66.02'
Set-TextValue 'E23' '  +1.84%  '
Set-TextValue 'E24' '  -0.32%  '
Set-TextValue 'E25' '  -4.89%  '
Set-TextValue 'D26' '572.34'
Set-TextValue 'E26' '  +1.63%  '
Set-TextValue 'E27' '  -0.03%  '
Set-TextValue 'D28' '2.530.36'
Set-TextValue 'E28' '  -0.45%  '
Set-TextValue 'D29' '0.0₃0943'
Set-TextValue 'E29' '  +1.03%  '
Set-TextValue 'E30' '  -2.17%  '
Set-TextValue 'D31' '1.42'
Set-TextValue 'E31' '  -2.47%  '
Set-TextValue 'E32' '  -0.04%  '
Set-TextValue 'E34' '  -2.30%  '
Set-TextValue 'D35' '1.00'
Set-TextValue 'E35' '  +0.11%  '
Set-TextValue 'E36' '  -2.67%  '
Set-TextValue 'D37' '5.50'
Set-TextValue 'E37' '  -4.67%  '
Set-TextValue 'B38' 'Monero'
Set-TextValue 'C38' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D38' '151.91'
Set-TextValue 'E38' '  +3.54%  '
Set-TextValue 'B39' 'PolygonEcosystemToken'
Set-TextValue 'C39' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D39' '0.379'
Set-TextValue 'E39' '  -1.32%  '
Set-TextValue 'D40' '18.63'
Set-TextValue 'E40' '  -0.77%  '
Set-TextValue 'E41' '  -9.81%  '
Set-TextValue 'D42' '0.993'
Set-TextValue 'E42' '  -0.78%  '
Set-TextValue 'E43' '  -1.47%  '
Set-TextValue 'D44' '148.05'
Set-TextValue 'E44' '  -1.86%  '
Set-TextValue 'D45' '3.63'
Set-TextValue 'E45' '  -0.18%  '
Set-TextValue 'D46' '0.0532'
Set-TextValue 'E46' '  -1.45%  '
Set-TextValue 'D47' '19.94'
Set-TextValue 'E47' '  -2.15%  '
Set-TextValue 'E48' '  +0.00%  '
Set-TextValue 'E49' '  +0.82%  '
Set-TextValue 'E50' '  -0.96%  '
Set-TextValue 'D51' '11.54'
Set-TextValue 'E51' '  +0.44%  '
